$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing content first (workbook had a different row layout before)
$ws.Cells.Clear()

# Write the updated localization key/value table
$ws.Cells.Item(1, 1).Value2 = 'keys'
$ws.Cells.Item(1, 2).Value2 = 'en'
$ws.Cells.Item(3, 1).Value2 = 'GAME_NAME'
$ws.Cells.Item(3, 2).Value2 = 'Danger Chaser'
$ws.Cells.Item(4, 1).Value2 = 'NEW_GAME'
$ws.Cells.Item(4, 2).Value2 = 'New Game'
$ws.Cells.Item(5, 1).Value2 = 'PLAY'
$ws.Cells.Item(5, 2).Value2 = 'Play'
$ws.Cells.Item(6, 1).Value2 = 'PRESS_ANY_BUTTON'
$ws.Cells.Item(6, 2).Value2 = 'Press any button'
$ws.Cells.Item(7, 1).Value2 = 'TUTORIAL'
$ws.Cells.Item(7, 2).Value2 = 'Tutorial'
$ws.Cells.Item(8, 1).Value2 = 'OPTIONS'
$ws.Cells.Item(8, 2).Value2 = 'Options'
$ws.Cells.Item(9, 1).Value2 = 'QUIT'
$ws.Cells.Item(9, 2).Value2 = 'Quit'
$ws.Cells.Item(10, 1).Value2 = 'LEVEL'
$ws.Cells.Item(10, 2).Value2 = 'Level'
$ws.Cells.Item(11, 1).Value2 = 'YES'
$ws.Cells.Item(11, 2).Value2 = 'Yes'
$ws.Cells.Item(12, 1).Value2 = 'NONE'
$ws.Cells.Item(12, 2).Value2 = 'No'
$ws.Cells.Item(13, 1).Value2 = 'RETURN'
$ws.Cells.Item(13, 2).Value2 = 'Return'
$ws.Cells.Item(14, 1).Value2 = 'PAUSED'
$ws.Cells.Item(14, 2).Value2 = 'Paused'
$ws.Cells.Item(15, 1).Value2 = 'RESUME'
$ws.Cells.Item(15, 2).Value2 = 'Resume'
$ws.Cells.Item(16, 1).Value2 = 'RESTART'
$ws.Cells.Item(16, 2).Value2 = 'Restart'
$ws.Cells.Item(17, 1).Value2 = 'ACTION_NOT_READY'
$ws.Cells.Item(17, 2).Value2 = 'Action not ready'
$ws.Cells.Item(18, 1).Value2 = 'INTERACTION_KEY'
$ws.Cells.Item(18, 2).Value2 = 'Z'
$ws.Cells.Item(19, 1).Value2 = 'MENU'
$ws.Cells.Item(19, 2).Value2 = 'Menu'
$ws.Cells.Item(20, 1).Value2 = 'MUSIC'
$ws.Cells.Item(20, 2).Value2 = 'Music'
$ws.Cells.Item(21, 1).Value2 = 'SFX'
$ws.Cells.Item(21, 2).Value2 = 'Sfx'
$ws.Cells.Item(22, 1).Value2 = 'AMBIENCE'
$ws.Cells.Item(22, 2).Value2 = 'Ambience'
$ws.Cells.Item(23, 1).Value2 = 'LOADING'
$ws.Cells.Item(23, 2).Value2 = 'Loading'
$ws.Cells.Item(24, 1).Value2 = 'CONTROLS'
$ws.Cells.Item(24, 2).Value2 = 'Controls'
$ws.Cells.Item(25, 1).Value2 = 'SETTINGS'
$ws.Cells.Item(25, 2).Value2 = 'Settings'
$ws.Cells.Item(26, 1).Value2 = 'BACK'
$ws.Cells.Item(26, 2).Value2 = 'Back'
$ws.Cells.Item(27, 1).Value2 = 'KEYBINDINGS'
$ws.Cells.Item(27, 2).Value2 = 'Keybindings'
$ws.Cells.Item(28, 1).Value2 = 'ACTION'
$ws.Cells.Item(28, 2).Value2 = 'Action'
$ws.Cells.Item(29, 1).Value2 = 'CONTROL'
$ws.Cells.Item(29, 2).Value2 = 'Key'
$ws.Cells.Item(30, 1).Value2 = 'CHANGE'
$ws.Cells.Item(30, 2).Value2 = 'Change'
$ws.Cells.Item(31, 1).Value2 = 'ui_up'
$ws.Cells.Item(31, 2).Value2 = 'Up'
$ws.Cells.Item(32, 1).Value2 = 'ui_left'
$ws.Cells.Item(32, 2).Value2 = 'Left'
$ws.Cells.Item(33, 1).Value2 = 'ui_right'
$ws.Cells.Item(33, 2).Value2 = 'Right'
$ws.Cells.Item(34, 1).Value2 = 'ui_down'
$ws.Cells.Item(34, 2).Value2 = 'Down'
$ws.Cells.Item(35, 1).Value2 = 'EXIT'
$ws.Cells.Item(35, 2).Value2 = 'Exit'
$ws.Cells.Item(36, 1).Value2 = 'CONTINUE'
$ws.Cells.Item(36, 2).Value2 = 'Continue'
$ws.Cells.Item(37, 1).Value2 = 'AUDIO'
$ws.Cells.Item(37, 2).Value2 = 'Audio'
$ws.Cells.Item(38, 1).Value2 = 'GRAPHICS'
$ws.Cells.Item(38, 2).Value2 = 'Graphics'
$ws.Cells.Item(39, 1).Value2 = 'SCREEN_SHAKE'
$ws.Cells.Item(39, 2).Value2 = 'Screen Shake'
$ws.Cells.Item(40, 1).Value2 = 'FRAME_FREEZE'
$ws.Cells.Item(40, 2).Value2 = 'Frame Freeze'
$ws.Cells.Item(41, 1).Value2 = 'ON'
$ws.Cells.Item(41, 2).Value2 = 'On'
$ws.Cells.Item(42, 1).Value2 = 'OFF'
$ws.Cells.Item(42, 2).Value2 = 'Off'
$ws.Cells.Item(43, 1).Value2 = 'DISABLED'
$ws.Cells.Item(43, 2).Value2 = 'Disabled'
$ws.Cells.Item(44, 1).Value2 = 'LOW'
$ws.Cells.Item(44, 2).Value2 = 'Low'
$ws.Cells.Item(45, 1).Value2 = 'NORMAL'
$ws.Cells.Item(45, 2).Value2 = 'Normal'
$ws.Cells.Item(46, 1).Value2 = 'HIGH'
$ws.Cells.Item(46, 2).Value2 = 'High'
$ws.Cells.Item(47, 1).Value2 = 'EXTREME'
$ws.Cells.Item(47, 2).Value2 = 'Extreme'
$ws.Cells.Item(48, 1).Value2 = 'VOMIT'
$ws.Cells.Item(48, 2).Value2 = 'Vomit'
$ws.Cells.Item(49, 1).Value2 = 'DECLINE'
$ws.Cells.Item(49, 2).Value2 = 'Decline'
$ws.Cells.Item(50, 1).Value2 = 'ACCEPT'
$ws.Cells.Item(50, 2).Value2 = 'Accept'
$ws.Cells.Item(100, 1).Value2 = 'QUIT_TO_MENU'
$ws.Cells.Item(100, 2).Value2 = 'Quit to Menu'
$ws.Cells.Item(101, 1).Value2 = 'QUESTS'
$ws.Cells.Item(101, 2).Value2 = 'Quests'
$ws.Cells.Item(102, 1).Value2 = 'ACT'
$ws.Cells.Item(102, 2).Value2 = 'Act'
$ws.Cells.Item(103, 1).Value2 = 'ZERO_LIMIT'
$ws.Cells.Item(103, 2).Value2 = 'Zero Limit'
$ws.Cells.Item(104, 1).Value2 = 'ZERO_LIMIT_ACT_1'
$ws.Cells.Item(104, 2).Value2 = 'Zero Limit Act 1'
$ws.Cells.Item(105, 1).Value2 = 'ZERO_LIMIT_ACT_1_DESCRIPTION'
$ws.Cells.Item(105, 2).Value2 = 'this is where I explain the story for now but nothing is here because reeeeeeeeee'
$ws.Cells.Item(106, 1).Value2 = 'ZERO_LIMIT_ACT_2'
$ws.Cells.Item(106, 2).Value2 = 'Zero Limit Act 2'
$ws.Cells.Item(107, 1).Value2 = 'ZERO_LIMIT_ACT_2_DESCRIPTION'
$ws.Cells.Item(107, 2).Value2 = 'insert description here lol'
$ws.Cells.Item(108, 1).Value2 = 'DEBUG_ENEMY'
$ws.Cells.Item(108, 2).Value2 = 'Debug Enemy'
$ws.Cells.Item(109, 1).Value2 = 'DEBUG_ENEMY_DESCRIPTION'
$ws.Cells.Item(109, 2).Value2 = 'I put enemies in here and test them out to make sure that they try to kill you the correct way.'

# Restore selection/viewport roughly matching the authored file
$ws.Range("A105").Select()

Write-Host "Localization strings sheet rebuilt."